$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 (ID 33): mark as complete
$ws.Range("D35").Value = $true

# Row 45 (ID 43): mark as complete, add reviewer note, row grows to fit wrapped note
$ws.Range("D45").Value = $true
$ws.Range("E45").Value = "better sign posting added and overviews added for non-technical readers"
$ws.Rows.Item(45).RowHeight = 28.8

# Row 51 (ID 49): mark as complete, add reviewer note referencing item 43
$ws.Range("D51").Value = $true
$ws.Range("E51").Value = "see item 43"

# Row 52 (ID 50): mark as complete
$ws.Range("D52").Value = $true

# Move the selection to reflect where the reviewer ended up working
$ws.Range("E52").Select()
